$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 591.05554
$ws.Range("I33").Value = 336.16666
$ws.Range("K33").Value = 336.16666
$ws.Range("M33").Value = -107.16666

$ws.Range("H43").Value = 2881.125
$ws.Range("I43").Value = 2046
$ws.Range("K43").Value = 2046
$ws.Range("M43").Value = -1977

$ws.Range("H88").Value = 1268
$ws.Range("J88").Value = 1652
$ws.Range("L88").Value = 1652
$ws.Range("N88").Value = -2464

$ws.Range("H91").Value = 1268
$ws.Range("J91").Value = 1652
$ws.Range("L91").Value = 1652
$ws.Range("N91").Value = -4460

$ws.Range("H94").Value = 1193.3334
$ws.Range("I94").Value = 790
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 790
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -339
$ws.Range("N94").Value = -2902

$ws.Range("H113").Value = 4741.702
$ws.Range("I113").Value = 4500.2563
$ws.Range("K113").Value = 4500.2563
$ws.Range("M113").Value = -1246.2563

$ws.Range("H127").Value = 825
$ws.Range("I127").Value = 450
$ws.Range("K127").Value = 1350
$ws.Range("M127").Value = 3610

$ws.Range("H129").Value = 1699
$ws.Range("J129").Value = 1699
$ws.Range("L129").Value = 5097
$ws.Range("N129").Value = -15097

$ws.Range("H132").Value = 1684.5238
$ws.Range("I132").Value = 1519
$ws.Range("K132").Value = 4557
$ws.Range("M132").Value = -2027

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 69000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H55").Value = 10048
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H61").Value = 3561.2334
$ws.Range("I61").Value = 3059.2666
$ws.Range("K61").Value = 3059.2666
$ws.Range("M61").Value = -2847.2666

$ws.Range("H132").Value = 2931.6553
$ws.Range("I132").Value = 2931.6553
$ws.Range("K132").Value = 8794.965899999999
$ws.Range("M132").Value = -6264.965899999999

$ws.Range("H136").Value = 3561.2334
$ws.Range("I136").Value = 3059.2666
$ws.Range("K136").Value = 9177.799800000001
$ws.Range("M136").Value = -6627.799800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2620.8
$ws.Range("I86").Value = 3941.5
$ws.Range("K86").Value = 3941.5
$ws.Range("M86").Value = -2818.5

$ws.Range("H89").Value = 2620.8
$ws.Range("I89").Value = 3941.5
$ws.Range("K89").Value = 19707.5
$ws.Range("M89").Value = -14091.5

$ws.Range("H134").Value = 3292.9
$ws.Range("I134").Value = 3297.6667
$ws.Range("K134").Value = 9893.000100000001
$ws.Range("M134").Value = -7358.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4725.9536
$ws.Range("I31").Value = 2512.1
$ws.Range("J31").Value = 5396.8184
$ws.Range("K31").Value = 2512.1
$ws.Range("L31").Value = 5396.8184
$ws.Range("M31").Value = -2217.1
$ws.Range("N31").Value = -5986.8184

$ws.Range("H34").Value = 4725.9536
$ws.Range("I34").Value = 2512.1
$ws.Range("J34").Value = 5396.8184
$ws.Range("K34").Value = 2512.1
$ws.Range("L34").Value = 5396.8184
$ws.Range("M34").Value = -2310.1
$ws.Range("N34").Value = -5800.8184

$ws.Range("H94").Value = 1680.2
$ws.Range("I94").Value = 1146.5
$ws.Range("K94").Value = 1146.5
$ws.Range("M94").Value = -695.5

$ws.Range("H132").Value = 4306
$ws.Range("I132").Value = 4306
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12918
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10388
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 470
$ws.Range("I33").Value = 216
$ws.Range("J33").Value = 724
$ws.Range("K33").Value = 1296
$ws.Range("L33").Value = 4344
$ws.Range("M33").Value = -1013
$ws.Range("N33").Value = -4910

$ws.Range("H113").Value = 1923.9048
$ws.Range("I113").Value = 915.7143
$ws.Range("J113").Value = 2428
$ws.Range("K113").Value = 2747.1429
$ws.Range("L113").Value = 7284
$ws.Range("M113").Value = -577.1428999999998
$ws.Range("N113").Value = -11624

$ws.Range("H115").Value = 2028.6666
$ws.Range("I115").Value = 2028
$ws.Range("J115").Value = 2029
$ws.Range("K115").Value = 6084
$ws.Range("L115").Value = 6087
$ws.Range("N115").Value = -8437
$ws.Range("M115").Value = -4909

$ws.Range("H131").Value = 1578.75
$ws.Range("J131").Value = 1678.1628
$ws.Range("L131").Value = 5034.4884
$ws.Range("N131").Value = -15114.4884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1216.4736
$ws.Range("I107").Value = 1028.75
$ws.Range("K107").Value = 1028.75
$ws.Range("M107").Value = 891.25

$ws.Range("H113").Value = 63199
$ws.Range("I113").Value = 55499.5
$ws.Range("J113").Value = 68332
$ws.Range("K113").Value = 55499.5
$ws.Range("L113").Value = 68332
$ws.Range("M113").Value = -53329.5
$ws.Range("N113").Value = -72672

$ws.Range("H126").Value = 2454.5356
$ws.Range("I126").Value = 1800.9375
$ws.Range("K126").Value = 5402.8125
$ws.Range("M126").Value = -2932.8125

$ws.Range("H132").Value = 3010.394
$ws.Range("I132").Value = 2652
$ws.Range("K132").Value = 7956
$ws.Range("M132").Value = -5426

$ws.Range("H133").Value = 65000
$ws.Range("J133").Value = 65000
$ws.Range("L133").Value = 65000
$ws.Range("N133").Value = -75120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2240.4092
$ws.Range("I16").Value = 2171.611
$ws.Range("K16").Value = 2171.611
$ws.Range("M16").Value = -2001.611

$ws.Range("H40").Value = 1496.5
$ws.Range("I40").Value = 1507.2222
$ws.Range("K40").Value = 1507.2222
$ws.Range("M40").Value = -1371.2222

$ws.Range("H55").Value = 1541.9375
$ws.Range("I55").Value = 1554.3684
$ws.Range("K55").Value = 1554.3684
$ws.Range("M55").Value = -1381.3684

$ws.Range("H61").Value = 3249.75
$ws.Range("I61").Value = 3166.3333
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 3166.3333
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -2964.3333
$ws.Range("N61").Value = -3904

$ws.Range("H113").Value = 3249.75
$ws.Range("I113").Value = 3166.3333
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 3166.3333
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -996.3332999999998
$ws.Range("N113").Value = -7840

$ws.Range("H122").Value = 14328.389
$ws.Range("I122").Value = 15307.615
$ws.Range("J122").Value = 11782.4
$ws.Range("K122").Value = 45922.845
$ws.Range("L122").Value = 35347.2
$ws.Range("M122").Value = -43472.845
$ws.Range("N122").Value = -40247.2

$ws.Range("H132").Value = 1201.3334
$ws.Range("I132").Value = 1201.3334
$ws.Range("K132").Value = 3604.0002
$ws.Range("M132").Value = -1074.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 120000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 120000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 120000
$ws.Range("N75").Value = -121872
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 120000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 120000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 360000
$ws.Range("N78").Value = -369360
$ws.Range("M78").ClearContents()

$ws.Range("H96").Value = 14094.5
$ws.Range("J96").Value = 14421.429
$ws.Range("L96").Value = 14421.429
$ws.Range("N96").Value = -17167.429

$ws.Range("H132").Value = 3180.1785
$ws.Range("I132").Value = 2668.5833
$ws.Range("J132").Value = 6249.75
$ws.Range("K132").Value = 8005.749899999999
$ws.Range("L132").Value = 18749.25
$ws.Range("M132").Value = -5475.749899999999
$ws.Range("N132").Value = -23809.25

$ws.Range("H136").Value = 2735.739
$ws.Range("I136").Value = 1977.6666
$ws.Range("J136").Value = 3223.0715
$ws.Range("K136").Value = 5932.9998
$ws.Range("L136").Value = 9669.2145
$ws.Range("M136").Value = -3382.9998
$ws.Range("N136").Value = -14769.2145
